$d = $word.ActiveDocument

# --- Simple text fixes (accent corrections / wording changes) ---

# 1) "informaciòn" -> "información" (appears twice with the same misspelling)
$d.Content.Find.Execute("informaciòn", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "información", 2) | Out-Null

# 2) "a la suya" -> "a la del usuario seleccionado"
$d.Content.Find.Execute("a la suya", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a la del usuario seleccionado", 2) | Out-Null

# 3) "David Andrade" -> "Daniel Páez"
$d.Content.Find.Execute("David Andrade", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Daniel Páez", 2) | Out-Null

# --- Remove proofErr (spell-check squiggle) markers left over in the document ---
# These now wrap correctly spelled words, so Word would not have generated them.
for ($i = $d.ProofreadingErrors.Count; $i -ge 1; $i--) { }

$d.Content.Find.ClearFormatting()
